$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (6th column) widens from 10.46875 to 11.71875 "characters".
# The COM ColumnWidth setter here only lands on multiples of 1/6 of a
# character, so 10.833333333333334 is the closest achievable value
# (serialises to the same 11.666... which is nearest to 11.71875).
$ws.Columns.Item(6).ColumnWidth = 10.833333333333334

# New 64-bit-design utilization numbers for row 2 (LUT, FF, BRAM, DSP).
$ws.Range("B2").Value = 50.75752258300781
$ws.Range("D2").Value = 19.862781524658203
$ws.Range("E2").Value = 57.85714340209961
$ws.Range("F2").Value = 72.7272720336914
